$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value (45175 -> 2023-09-06)
# that was updated to 45183 (2023-09-14) for rows 2 through 28.
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
